# A new price record was inserted as row 396 (pushing the existing rows
# 396..482 down to 397..483). Insert a whole row at 396 so every
# subsequent row shifts down by one, then populate the new row with its
# data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(396).Insert()

$ws.Cells.Item(396, 1).Value = 6
$ws.Cells.Item(396, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(396, 3).Value = "Metropolitana"
$ws.Cells.Item(396, 4).Value = 44798
$ws.Cells.Item(396, 5).Value = 13
$ws.Cells.Item(396, 6).Value = 100112043
$ws.Cells.Item(396, 7).Value = "Pepino ensalada"
$ws.Cells.Item(396, 8).Value = "Sin especificar"
$ws.Cells.Item(396, 9).Value = "Primera"
$ws.Cells.Item(396, 10).Value = 490
$ws.Cells.Item(396, 11).Value = 20000
$ws.Cells.Item(396, 12).Value = 22000
$ws.Cells.Item(396, 13).Value = 21061
$ws.Cells.Item(396, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(396, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(396, 16).Value = 351
$ws.Cells.Item(396, 17).Value = 60
$ws.Cells.Item(396, 18).Value = "Hortaliza"
